# Generate Report for Handoff
# Adds two new source files (85ee349b... and bff880de...) to the
# localization-status workbook, pushing ".localization-config" down to
# row 6 on every sheet, and records their "Ready for handoff" /
# "Include" handoff rows on the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$baseMd   = "https://github.com/OpenLocalizationTest/oltest/blob/9f07fe290603636178a7c845719e760545101205/e2e/"
$baseCfg  = "https://github.com/OpenLocalizationTest/oltest/blob/9f07fe290603636178a7c845719e760545101205/.localization-config"
$baseZh   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c539178c125b7d9348b005eb7fac4e3172ceb06e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/"
$baseDe   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea1279ed0a828b435eff7a3b131c8ba2de067b3c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/"

$guid1 = "85ee349b-2501-4ac1-82f8-284c63e1198b"
$guid2 = "bff880de-8f5b-4d1d-8aed-8b68bfdd427e"
$hash1 = "3c945d6cc0be654d3baf95d2c7ea878da09d5fee"
$hash2 = "5a014e88b6ff33c56dfdbd1fbdcc893edc5b4e0a"

$md1 = "$guid1.md"
$md2 = "$guid2.md"
$cfg = ".localization-config"

$xlf1zh = "$guid1.$hash1.zh-cn.xlf"
$xlf2zh = "$guid2.$hash2.zh-cn.xlf"
$xlf1de = "$guid1.$hash1.de-de.xlf"
$xlf2de = "$guid2.$hash2.de-de.xlf"

$dt1 = "2016-02-26 05:22:25"
$dt2 = "2016-02-26 05:22:36"
$never = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A4").Value = $md1
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = $md2
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"

$ws1.Range("A6").Value = $cfg
$ws1.Range("B6").Value = "Not to be localized"
$ws1.Range("C6").Value = "Not to be localized"

# Rebuild hyperlinks (the engine's Hyperlinks.Delete() clears the whole
# sheet collection, so all links - old and new - are re-added in order).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $baseMd + "1b5d57df-ea1d-4c26-8ccd-6e31db159268.md", "", "", "1b5d57df-ea1d-4c26-8ccd-6e31db159268.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $baseMd + "3c1d0466-a578-406f-ae57-5e2575653435.md", "", "", "3c1d0466-a578-406f-ae57-5e2575653435.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), $baseMd + $md1, "", "", $md1)
$ws1.Hyperlinks.Add($ws1.Range("A5"), $baseMd + $md2, "", "", $md2)
$ws1.Hyperlinks.Add($ws1.Range("A6"), $baseCfg, "", "", $cfg)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A4").Value = $md1
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = $xlf1zh
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D4").Value = $dt1
$ws2.Range("G4").Value = $never
$ws2.Range("H4").Value = "Include"

$ws2.Range("A5").Value = $md2
$ws2.Range("B5").Value = "Ready for handoff"
$ws2.Range("C5").Value = $xlf2zh
$ws2.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D5").Value = $dt1
$ws2.Range("G5").Value = $never
$ws2.Range("H5").Value = "Include"

$ws2.Range("A6").Value = $cfg
$ws2.Range("B6").Value = "Not to be localized"
$ws2.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D6").Value = $never
$ws2.Range("G6").Value = $never
$ws2.Range("H6").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $baseMd + "1b5d57df-ea1d-4c26-8ccd-6e31db159268.md", "", "", "1b5d57df-ea1d-4c26-8ccd-6e31db159268.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), $baseZh + "1b5d57df-ea1d-4c26-8ccd-6e31db159268.22d7d7ced49f4033dd3add9bd40746d4366b2e3c.zh-cn.xlf", "", "", "1b5d57df-ea1d-4c26-8ccd-6e31db159268.22d7d7ced49f4033dd3add9bd40746d4366b2e3c.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $baseMd + "3c1d0466-a578-406f-ae57-5e2575653435.md", "", "", "3c1d0466-a578-406f-ae57-5e2575653435.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), $baseZh + "3c1d0466-a578-406f-ae57-5e2575653435.39999fc6db05c4892588f7cc4e1c31417b50fd05.zh-cn.xlf", "", "", "3c1d0466-a578-406f-ae57-5e2575653435.39999fc6db05c4892588f7cc4e1c31417b50fd05.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), $baseMd + $md1, "", "", $md1)
$ws2.Hyperlinks.Add($ws2.Range("C4"), $baseZh + $xlf1zh, "", "", $xlf1zh)
$ws2.Hyperlinks.Add($ws2.Range("A5"), $baseMd + $md2, "", "", $md2)
$ws2.Hyperlinks.Add($ws2.Range("C5"), $baseZh + $xlf2zh, "", "", $xlf2zh)
$ws2.Hyperlinks.Add($ws2.Range("A6"), $baseCfg, "", "", $cfg)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A4").Value = $md1
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = $xlf1de
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D4").Value = $dt2
$ws3.Range("G4").Value = $never
$ws3.Range("H4").Value = "Include"

$ws3.Range("A5").Value = $md2
$ws3.Range("B5").Value = "Ready for handoff"
$ws3.Range("C5").Value = $xlf2de
$ws3.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D5").Value = $dt2
$ws3.Range("G5").Value = $never
$ws3.Range("H5").Value = "Include"

$ws3.Range("A6").Value = $cfg
$ws3.Range("B6").Value = "Not to be localized"
$ws3.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D6").Value = $never
$ws3.Range("G6").Value = $never
$ws3.Range("H6").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $baseMd + "1b5d57df-ea1d-4c26-8ccd-6e31db159268.md", "", "", "1b5d57df-ea1d-4c26-8ccd-6e31db159268.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), $baseDe + "1b5d57df-ea1d-4c26-8ccd-6e31db159268.22d7d7ced49f4033dd3add9bd40746d4366b2e3c.de-de.xlf", "", "", "1b5d57df-ea1d-4c26-8ccd-6e31db159268.22d7d7ced49f4033dd3add9bd40746d4366b2e3c.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $baseMd + "3c1d0466-a578-406f-ae57-5e2575653435.md", "", "", "3c1d0466-a578-406f-ae57-5e2575653435.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), $baseDe + "3c1d0466-a578-406f-ae57-5e2575653435.39999fc6db05c4892588f7cc4e1c31417b50fd05.de-de.xlf", "", "", "3c1d0466-a578-406f-ae57-5e2575653435.39999fc6db05c4892588f7cc4e1c31417b50fd05.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), $baseMd + $md1, "", "", $md1)
$ws3.Hyperlinks.Add($ws3.Range("C4"), $baseDe + $xlf1de, "", "", $xlf1de)
$ws3.Hyperlinks.Add($ws3.Range("A5"), $baseMd + $md2, "", "", $md2)
$ws3.Hyperlinks.Add($ws3.Range("C5"), $baseDe + $xlf2de, "", "", $xlf2de)
$ws3.Hyperlinks.Add($ws3.Range("A6"), $baseCfg, "", "", $cfg)
